$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.044767888252353
$ws.Cells.Item(2, 4).Value = 1.054120968337501
$ws.Cells.Item(2, 5).Value = 1.04257931755384
$ws.Cells.Item(2, 6).Value = 1.061064280785392
$ws.Cells.Item(2, 9).Value = 1.037446097882198
$ws.Cells.Item(2, 10).Value = 1.049831207851448
$ws.Cells.Item(2, 11).Value = 1.056865052911502
$ws.Cells.Item(2, 12).Value = 1.045355606286391
$ws.Cells.Item(2, 13).Value = 1.063789357252356
$ws.Cells.Item(2, 14).Value = 1.020393081036093
# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.046151503999267
$ws.Cells.Item(3, 4).Value = 1.055478553076081
$ws.Cells.Item(3, 5).Value = 1.043767602474577
$ws.Cells.Item(3, 6).Value = 1.062572484014418
$ws.Cells.Item(3, 9).Value = 1.037736802806428
$ws.Cells.Item(3, 10).Value = 1.05086006696713
$ws.Cells.Item(3, 11).Value = 1.058034426019085
$ws.Cells.Item(3, 12).Value = 1.046353763177557
$ws.Cells.Item(3, 13).Value = 1.065110360046207
$ws.Cells.Item(3, 14).Value = 1.020740899390178
# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.047045332084307
$ws.Cells.Item(4, 4).Value = 1.056355873264885
$ws.Cells.Item(4, 5).Value = 1.044535552590347
$ws.Cells.Item(4, 6).Value = 1.063547509923694
$ws.Cells.Item(4, 9).Value = 1.03792229435905
$ws.Cells.Item(4, 10).Value = 1.051523908247096
$ws.Cells.Item(4, 11).Value = 1.058789427034369
$ws.Cells.Item(4, 12).Value = 1.046998121678641
$ws.Cells.Item(4, 13).Value = 1.065963749451744
$ws.Cells.Item(4, 14).Value = 1.020965153400273
# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.0474207526779
$ws.Cells.Item(5, 4).Value = 1.056724433694723
$ws.Cells.Item(5, 5).Value = 1.044858175495398
$ws.Cells.Item(5, 6).Value = 1.063957205358897
$ws.Cells.Item(5, 9).Value = 1.037999650508022
$ws.Cells.Item(5, 10).Value = 1.051802536453063
$ws.Cells.Item(5, 11).Value = 1.059106436224544
$ws.Cells.Item(5, 12).Value = 1.047268650847432
$ws.Cells.Item(5, 13).Value = 1.066322187288242
$ws.Cells.Item(5, 14).Value = 1.021059237722066
# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.047483767397329
$ws.Cells.Item(6, 4).Value = 1.056786301220221
$ws.Cells.Item(6, 5).Value = 1.044912332321571
$ws.Cells.Item(6, 6).Value = 1.064025983144893
$ws.Cells.Item(6, 9).Value = 1.038012602361941
$ws.Cells.Item(6, 10).Value = 1.051849293044342
$ws.Cells.Item(6, 11).Value = 1.059159640534407
$ws.Cells.Item(6, 12).Value = 1.047314052932614
$ws.Cells.Item(6, 13).Value = 1.066382351547833
$ws.Cells.Item(6, 14).Value = 1.021075023668723
# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.047050349820176
$ws.Cells.Item(7, 4).Value = 1.056360799022344
$ws.Cells.Item(7, 5).Value = 1.044539864365936
$ws.Cells.Item(7, 6).Value = 1.063552985094298
$ws.Cells.Item(7, 9).Value = 1.037923330447997
$ws.Cells.Item(7, 10).Value = 1.051527633056672
$ws.Cells.Item(7, 11).Value = 1.058793664466239
$ws.Cells.Item(7, 12).Value = 1.047001737909245
$ws.Cells.Item(7, 13).Value = 1.065968540194464
$ws.Cells.Item(7, 14).Value = 1.020966411312659
# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.045235794213956
$ws.Cells.Item(8, 4).Value = 1.054580007047942
$ws.Cells.Item(8, 5).Value = 1.042981103082239
$ws.Cells.Item(8, 6).Value = 1.06157417134279
$ws.Cells.Item(8, 9).Value = 1.037544884995931
$ws.Cells.Item(8, 10).Value = 1.050179311327046
$ws.Cells.Item(8, 11).Value = 1.057260594817613
$ws.Cells.Item(8, 12).Value = 1.045693254181241
$ws.Cells.Item(8, 13).Value = 1.064236086786222
$ws.Cells.Item(8, 14).Value = 1.020510795911273
# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.042026847120366
$ws.Cells.Item(9, 4).Value = 1.051433157993148
$ws.Cells.Item(9, 5).Value = 1.040226911068205
$ws.Cells.Item(9, 6).Value = 1.05808022916023
$ws.Cells.Item(9, 9).Value = 1.036857936874992
$ws.Cells.Item(9, 10).Value = 1.047788666714661
$ws.Cells.Item(9, 11).Value = 1.054546199300147
$ws.Cells.Item(9, 12).Value = 1.043375761493396
$ws.Cells.Item(9, 13).Value = 1.061172421692568
$ws.Cells.Item(9, 14).Value = 1.019701698169574
# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.03987947584111
$ws.Cells.Item(10, 4).Value = 1.049328976094148
$ws.Cells.Item(10, 5).Value = 1.038385512790063
$ws.Cells.Item(10, 6).Value = 1.055745846531588
$ws.Cells.Item(10, 9).Value = 1.036386386563624
$ws.Cells.Item(10, 10).Value = 1.046184741939964
$ws.Cells.Item(10, 11).Value = 1.052727618134152
$ws.Cells.Item(10, 12).Value = 1.041822609804244
$ws.Cells.Item(10, 13).Value = 1.059122347887987
$ws.Cells.Item(10, 14).Value = 1.019158017134528
# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.038947643925342
$ws.Cells.Item(11, 4).Value = 1.048416279463728
$ws.Cells.Item(11, 5).Value = 1.037586858709111
$ws.Cells.Item(11, 6).Value = 1.054733743433724
$ws.Cells.Item(11, 9).Value = 1.036178958232003
$ws.Cells.Item(11, 10).Value = 1.045487756134126
$ws.Cells.Item(11, 11).Value = 1.051937955977995
$ws.Cells.Item(11, 12).Value = 1.041148089466714
$ws.Cells.Item(11, 13).Value = 1.058232759494989
$ws.Cells.Item(11, 14).Value = 1.018921562530947
# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.038601211549582
$ws.Cells.Item(12, 4).Value = 1.048077020768049
$ws.Cells.Item(12, 5).Value = 1.037290000114343
$ws.Cells.Item(12, 6).Value = 1.054357600721744
$ws.Cells.Item(12, 9).Value = 1.036101420953783
$ws.Cells.Item(12, 10).Value = 1.045228487463503
$ws.Cells.Item(12, 11).Value = 1.051644303333829
$ws.Cells.Item(12, 12).Value = 1.040897238048394
$ws.Cells.Item(12, 13).Value = 1.057902035885969
$ws.Cells.Item(12, 14).Value = 1.018833575270203
# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.038675536546231
$ws.Cells.Item(13, 4).Value = 1.048149803987537
$ws.Cells.Item(13, 5).Value = 1.037353686555459
$ws.Cells.Item(13, 6).Value = 1.054438293893438
$ws.Cells.Item(13, 9).Value = 1.036118075127798
$ws.Cells.Item(13, 10).Value = 1.045284118621626
$ws.Cells.Item(13, 11).Value = 1.051707308163953
$ws.Cells.Item(13, 12).Value = 1.040951060386518
$ws.Cells.Item(13, 13).Value = 1.057972990489985
$ws.Cells.Item(13, 14).Value = 1.018852455990273
# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.038919014052146
$ws.Cells.Item(14, 4).Value = 1.048388241218606
$ws.Cells.Item(14, 5).Value = 1.037562324454151
$ws.Cells.Item(14, 6).Value = 1.054702655534051
$ws.Cells.Item(14, 9).Value = 1.036172558965596
$ws.Cells.Item(14, 10).Value = 1.045466332627802
$ws.Cells.Item(14, 11).Value = 1.051913689466186
$ws.Cells.Item(14, 12).Value = 1.041127360250184
$ws.Cells.Item(14, 13).Value = 1.058205427749679
$ws.Cells.Item(14, 14).Value = 1.018914292699084
# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.039068987483913
$ws.Cells.Item(15, 4).Value = 1.04853511793748
$ws.Cells.Item(15, 5).Value = 1.03769084610796
$ws.Cells.Item(15, 6).Value = 1.054865510407257
$ws.Cells.Item(15, 9).Value = 1.036206063383633
$ws.Cells.Item(15, 10).Value = 1.045578550565355
$ws.Cells.Item(15, 11).Value = 1.052040802974396
$ws.Cells.Item(15, 12).Value = 1.041235943918692
$ws.Cells.Item(15, 13).Value = 1.058348601278647
$ws.Cells.Item(15, 14).Value = 1.018952371416515
# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.039941275584954
$ws.Cells.Item(16, 4).Value = 1.049389515102448
$ws.Cells.Item(16, 5).Value = 1.038438488658251
$ws.Cells.Item(16, 6).Value = 1.055812988474075
$ws.Cells.Item(16, 9).Value = 1.036400084400199
$ws.Cells.Item(16, 10).Value = 1.046230945991911
$ws.Cells.Item(16, 11).Value = 1.052779978474076
$ws.Cells.Item(16, 12).Value = 1.041867333052662
$ws.Cells.Item(16, 13).Value = 1.059181346506387
$ws.Cells.Item(16, 14).Value = 1.019173687854024
# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.040487896787162
$ws.Cells.Item(17, 4).Value = 1.049925030536201
$ws.Cells.Item(17, 5).Value = 1.038907109017116
$ws.Cells.Item(17, 6).Value = 1.056406963286143
$ws.Cells.Item(17, 9).Value = 1.036520918904818
$ws.Cells.Item(17, 10).Value = 1.046639509946372
$ws.Cells.Item(17, 11).Value = 1.053243049477976
$ws.Cells.Item(17, 12).Value = 1.042262849212082
$ws.Cells.Item(17, 13).Value = 1.059703194242266
$ws.Cells.Item(17, 14).Value = 1.019312234961071
# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.040806538624165
$ws.Cells.Item(18, 4).Value = 1.050237236388173
$ws.Cells.Item(18, 5).Value = 1.039180320672679
$ws.Cells.Item(18, 6).Value = 1.056753293805572
$ws.Cells.Item(18, 9).Value = 1.036591086757535
$ws.Cells.Item(18, 10).Value = 1.046877579861572
$ws.Cells.Item(18, 11).Value = 1.053512938649149
$ws.Cells.Item(18, 12).Value = 1.042493354891995
$ws.Cells.Item(18, 13).Value = 1.06000739700125
$ws.Cells.Item(18, 14).Value = 1.019392947191139
# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.040915154724091
$ws.Cells.Item(19, 4).Value = 1.050343665053671
$ws.Cells.Item(19, 5).Value = 1.039273457455072
$ws.Cells.Item(19, 6).Value = 1.056871362588746
$ws.Cells.Item(19, 9).Value = 1.036614959167632
$ws.Cells.Item(19, 10).Value = 1.046958715246746
$ws.Cells.Item(19, 11).Value = 1.05360492804337
$ws.Cells.Item(19, 12).Value = 1.042571918913417
$ws.Cells.Item(19, 13).Value = 1.060111091572004
$ws.Cells.Item(19, 14).Value = 1.019420451072309
# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.040429269526424
$ws.Cells.Item(20, 4).Value = 1.049867590487852
$ws.Cells.Item(20, 5).Value = 1.038856843626998
$ws.Cells.Item(20, 6).Value = 1.05634324838574
$ws.Cells.Item(20, 9).Value = 1.036507986880627
$ws.Cells.Item(20, 10).Value = 1.046595699610005
$ws.Cells.Item(20, 11).Value = 1.053193388343798
$ws.Cells.Item(20, 12).Value = 1.04222043398633
$ws.Cells.Item(20, 13).Value = 1.059647223811388
$ws.Cells.Item(20, 14).Value = 1.019297380518506
# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.038847324569234
$ws.Cells.Item(21, 4).Value = 1.048318034143274
$ws.Cells.Item(21, 5).Value = 1.037500891435814
$ws.Cells.Item(21, 6).Value = 1.054624813269998
$ws.Cells.Item(21, 9).Value = 1.036156528351101
$ws.Cells.Item(21, 10).Value = 1.045412685596196
$ws.Cells.Item(21, 11).Value = 1.051852924664827
$ws.Cells.Item(21, 12).Value = 1.041075452779388
$ws.Cells.Item(21, 13).Value = 1.058136988851853
$ws.Cells.Item(21, 14).Value = 1.018896087690689
# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.037850904451635
$ws.Cells.Item(22, 4).Value = 1.04734236017862
$ws.Cells.Item(22, 5).Value = 1.036647173934504
$ws.Cells.Item(22, 6).Value = 1.053543189340241
$ws.Cells.Item(22, 9).Value = 1.035932721386371
$ws.Cells.Item(22, 10).Value = 1.044666693309787
$ws.Cells.Item(22, 11).Value = 1.051008169512937
$ws.Cells.Item(22, 12).Value = 1.040353793340506
$ws.Cells.Item(22, 13).Value = 1.057185758050193
$ws.Cells.Item(22, 14).Value = 1.018642866916112
# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.038379297538676
$ws.Cells.Item(23, 4).Value = 1.047859718853121
$ws.Cells.Item(23, 5).Value = 1.037099858979443
$ws.Cells.Item(23, 6).Value = 1.054116692514405
$ws.Cells.Item(23, 9).Value = 1.036051634648554
$ws.Cells.Item(23, 10).Value = 1.045062366724454
$ws.Cells.Item(23, 11).Value = 1.051456177205308
$ws.Cells.Item(23, 12).Value = 1.040736527494348
$ws.Cells.Item(23, 13).Value = 1.057690185501072
$ws.Cells.Item(23, 14).Value = 1.018777191054674
# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.040455761258834
$ws.Cells.Item(24, 4).Value = 1.04989354564089
$ws.Cells.Item(24, 5).Value = 1.038879556784883
$ws.Cells.Item(24, 6).Value = 1.056372038793062
$ws.Cells.Item(24, 9).Value = 1.036513831272598
$ws.Cells.Item(24, 10).Value = 1.04661549635253
$ws.Cells.Item(24, 11).Value = 1.053215828729026
$ws.Cells.Item(24, 12).Value = 1.042239600195282
$ws.Cells.Item(24, 13).Value = 1.059672515003182
$ws.Cells.Item(24, 14).Value = 1.019304092910657
# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.042857833056588
$ws.Cells.Item(25, 4).Value = 1.052247777500198
$ws.Cells.Item(25, 5).Value = 1.04093984560758
$ws.Cells.Item(25, 6).Value = 1.058984366544906
$ws.Cells.Item(25, 9).Value = 1.037037917539802
$ws.Cells.Item(25, 10).Value = 1.048408477618086
$ws.Cells.Item(25, 11).Value = 1.055249497689792
$ws.Cells.Item(25, 12).Value = 1.043976308947873
$ws.Cells.Item(25, 13).Value = 1.061965773232145
$ws.Cells.Item(25, 14).Value = 1.019911617976991

Write-Host "Updated vm_pu values for 380 kV case"
